$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.256674408912659
$ws.Range("B1").Value = 2.593246698379517
$ws.Range("C1").Value = 4.996026992797852
$ws.Range("D1").Value = 1.998872756958008
$ws.Range("E1").Value = 1.156558632850647
